# Insert a new "Match ID" column at the front of the sheet, shifting all
# existing data one column to the right, then populate the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; this shifts data + styles
# (including the merged header cells) one column to the right automatically.
$ws.Range("A1").EntireColumn.Insert()

# Row 3 holds the (visible) column header labels - add the new header.
$ws.Range("A3").Value = "Match ID"

# Data rows 4-18 (one per player) all reference match id "2" for this
# particular game sheet.
$ws.Range("A4:A18").Value = 2

# Row 19 is the hidden totals row; temporarily unhide it while writing so
# the engine doesn't stamp a recalculated row height onto a hidden row.
$ws.Rows(19).Hidden = $false
$ws.Range("A19").Value = 2
$ws.Rows(19).Hidden = $true

# Header + data cells get the bold "label" styling used elsewhere in the
# sheet (bold font, no border), matching the rest of row 3 onward; the
# totals row (19) keeps the plain default style (untouched).
$ws.Range("A3:A18").Font.Bold = $true
$ws.Range("A3:A18").Borders.LineStyle = -4142

# Restore the selection to the newly added Match ID column's data range.
$ws.Range("A3:A18").Select()
